# Spring 23 week 14 inputs: append 23 new matchup rows to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Player_1, B=Points_1, C=Player_2, D=Points_2
$data = @(
    @(5,2,5,0),
    @(3,0,3,3),
    @(5,1,5,2),
    @(6,0,6,3),
    @(4,1,3,2),
    @(5,1,4,2),
    @(6,3,5,0),
    @(4,3,5,0),
    @(3,1,4,2),
    @(4,2,4,0),
    @(7,1,5,2),
    @(2,0,3,3),
    @(7,2,5,1),
    @(2,1,5,2),
    @(6,2,6,0),
    @(3,2,4,1),
    @(7,1,6,2),
    @(5,2,4,1),
    @(3,0,3,3),
    @(3,1,3,2),
    @(4,1,4,2),
    @(7,2,6,0),
    @(5,2,5,1)
)

$startRow = 1853
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Move the active selection to the next empty row, same as Excel does after
# typing the last row of data, and scroll the window to follow it.
$nextRow = $endRow + 1
$excel.ActiveWindow.ScrollRow = 1860
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A$nextRow").Select()
